# Insert a new row at row 114 (pushing existing rows 114:217 down to 115:218)
# and populate it with the new weekly data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(114).Insert()

$ws.Cells.Item(114, 1).Value = 6
$ws.Cells.Item(114, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(114, 3).Value = "Metropolitana"
$ws.Cells.Item(114, 4).Value = 44740
$ws.Cells.Item(114, 5).Value = 13
$ws.Cells.Item(114, 6).Value = 100112001
$ws.Cells.Item(114, 7).Value = "Berenjena"
$ws.Cells.Item(114, 8).Value = "Sin especificar"
$ws.Cells.Item(114, 9).Value = "Primera"
$ws.Cells.Item(114, 10).Value = 400
$ws.Cells.Item(114, 11).Value = 10000
$ws.Cells.Item(114, 12).Value = 12000
$ws.Cells.Item(114, 13).Value = 10850
$ws.Cells.Item(114, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(114, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(114, 16).Value = 217
$ws.Cells.Item(114, 17).Value = 50
$ws.Cells.Item(114, 18).Value = "Hortaliza"
